$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B1").Value = "id_direcao_saida"
$ws.Range("C1").Value = "direcao_saida"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Data rows: A = new id value, B = sequential index (was A), C = text (was B)
$data = @(
    @(43409, 1, "Carga Piezométrica"),
    @(48808, 2, "Cota Piezométrica"),
    @(7271, 3, "Cota do NA"),
    @(23375, 4, "Deslocamento Longitudinal (X)"),
    @(22596, 5, "Deslocamento Transversal (Y)"),
    @(45750, 6, "Deslocamento da Estaca Dir(-)/Esq(+)"),
    @(20547, 7, "Deslocamento do Afastamento Mont(-)/Jus(+)"),
    @(28397, 8, "Leitura"),
    @(33011, 9, "Recalque"),
    @(23084, 10, "Recalque (Z)"),
    @(33897, 11, "Vazão")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# Remove the now-unused row 13 (previous data had 12 rows, now has 11)
$ws.Rows.Item(13).Delete()
